$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 5, pushing existing rows 5-11 down to 6-12.
# Excel's Insert copies the formatting of the row above (row 4) by default.
$ws.Rows.Item(5).Insert()

# Set the text for the new row's first cell.
$ws.Range("A5").Value = "{#d = d.date}"

# Move the active selection as in the authored workbook.
$ws.Range("A15").Select()
